$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2

$old1 = "✅ 1000 Bs = 1.88 = 6868.23 pesos`n✅ 6868.23 pesos = 1.87 = 916.76 Bs"
$new1 = "✅ 1000 Bs = 1.93 = 7079.34 pesos`n✅ 7079.34 pesos = 1.92 = 965.75 Bs"
$text = $text.Replace($old1, $new1)

$cell.Value2 = $text

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value2 = 516.999
$ws2.Range("O10").Value2 = 3660.01
$ws2.Range("N12").Value2 = 3685
$ws2.Range("O12").Value2 = 502.7
